$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 6283.3335  # H40: 3750 -> 6283.3335
$ws.Cells.Item(40, 9).Value = 4300  # I40: 2899.8333 -> 4300
$ws.Cells.Item(40, 10).Value = 8266.666999999999  # J40: 4600.1665 -> 8266.666999999999
$ws.Cells.Item(40, 11).Value = 4300  # K40: 2899.8333 -> 4300
$ws.Cells.Item(40, 12).Value = 8266.666999999999  # L40: 4600.1665 -> 8266.666999999999
$ws.Cells.Item(40, 13).Value = -4125  # M40: -2724.8333 -> -4125
$ws.Cells.Item(40, 14).Value = -8616.666999999999  # N40: -4950.1665 -> -8616.666999999999
$ws.Cells.Item(86, 8).Value = 3098498.2  # H86: 3511535.8 -> 3098498.2
$ws.Cells.Item(86, 9).Value = 2012.4286  # I86: 2800 -> 2012.4286
$ws.Cells.Item(86, 10).Value = 5266038.5  # J86: 4787439.5 -> 5266038.5
$ws.Cells.Item(86, 11).Value = 2012.4286  # K86: 2800 -> 2012.4286
$ws.Cells.Item(86, 12).Value = 5266038.5  # L86: 4787439.5 -> 5266038.5
$ws.Cells.Item(86, 13).Value = -889.4286  # M86: -1677 -> -889.4286
$ws.Cells.Item(86, 14).Value = -5268284.5  # N86: -4789685.5 -> -5268284.5
$ws.Cells.Item(89, 8).Value = 3098498.2  # H89: 3511535.8 -> 3098498.2
$ws.Cells.Item(89, 9).Value = 2012.4286  # I89: 2800 -> 2012.4286
$ws.Cells.Item(89, 10).Value = 5266038.5  # J89: 4787439.5 -> 5266038.5
$ws.Cells.Item(89, 11).Value = 10062.143  # K89: 14000 -> 10062.143
$ws.Cells.Item(89, 12).Value = 26330192.5  # L89: 23937197.5 -> 26330192.5
$ws.Cells.Item(89, 13).Value = -4446.143  # M89: -8384 -> -4446.143
$ws.Cells.Item(89, 14).Value = -26341424.5  # N89: -23948429.5 -> -26341424.5
$ws.Cells.Item(112, 8).Value = 1332.7297  # H112: 1340.697 -> 1332.7297
$ws.Cells.Item(112, 10).Value = 1311.4166  # J112: 1316.9688 -> 1311.4166
$ws.Cells.Item(112, 12).Value = 3934.2498  # L112: 3950.9064 -> 3934.2498
$ws.Cells.Item(112, 14).Value = -6150.2498  # N112: -6166.9064 -> -6150.2498

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15877296  # H32: 15629218 -> 15877296
$ws.Cells.Item(32, 9).Value = 16670911  # I32: 16397622 -> 16670911
$ws.Cells.Item(32, 11).Value = 16670911  # K32: 16397622 -> 16670911
$ws.Cells.Item(32, 13).Value = -16670624  # M32: -16397335 -> -16670624
$ws.Cells.Item(45, 8).Value = 7291.1  # H45: 7842.4707 -> 7291.1
$ws.Cells.Item(45, 9).Value = 4560.7144  # I45: 4904.1665 -> 4560.7144
$ws.Cells.Item(45, 10).Value = 8761.308000000001  # J45: 9445.182000000001 -> 8761.308000000001
$ws.Cells.Item(45, 11).Value = 4560.7144  # K45: 4904.1665 -> 4560.7144
$ws.Cells.Item(45, 12).Value = 8761.308000000001  # L45: 9445.182000000001 -> 8761.308000000001
$ws.Cells.Item(45, 13).Value = -4183.7144  # M45: -4527.1665 -> -4183.7144
$ws.Cells.Item(45, 14).Value = -9515.308000000001  # N45: -10199.182 -> -9515.308000000001
$ws.Cells.Item(61, 8).Value = 4571.433  # H61: 4010.3845 -> 4571.433
$ws.Cells.Item(61, 9).Value = 1922.25  # I61: 1843.3334 -> 1922.25
$ws.Cells.Item(61, 10).Value = 5534.773  # J61: 8886.25 -> 5534.773
$ws.Cells.Item(61, 11).Value = 1922.25  # K61: 1843.3334 -> 1922.25
$ws.Cells.Item(61, 12).Value = 5534.773  # L61: 8886.25 -> 5534.773
$ws.Cells.Item(61, 13).Value = -1710.25  # M61: -1631.3334 -> -1710.25
$ws.Cells.Item(61, 14).Value = -5958.773  # N61: -9310.25 -> -5958.773
$ws.Cells.Item(74, 8).Value = 2506.875  # H74: 2508.25 -> 2506.875
$ws.Cells.Item(74, 9).Value = 2150.7144  # I74: 2152.2856 -> 2150.7144
$ws.Cells.Item(74, 11).Value = 2150.7144  # K74: 2152.2856 -> 2150.7144
$ws.Cells.Item(74, 13).Value = -1276.7144  # M74: -1278.2856 -> -1276.7144
$ws.Cells.Item(77, 8).Value = 2506.875  # H77: 2508.25 -> 2506.875
$ws.Cells.Item(77, 9).Value = 2150.7144  # I77: 2152.2856 -> 2150.7144
$ws.Cells.Item(77, 11).Value = 10753.572  # K77: 10761.428 -> 10753.572
$ws.Cells.Item(77, 13).Value = -6385.572  # M77: -6393.428 -> -6385.572
$ws.Cells.Item(92, 8).Value = 52241  # H92: 56684.168 -> 52241
$ws.Cells.Item(92, 10).Value = 52241  # J92: 56684.168 -> 52241
$ws.Cells.Item(92, 12).Value = 52241  # L92: 56684.168 -> 52241
$ws.Cells.Item(92, 14).Value = -57233  # N92: -61676.168 -> -57233
$ws.Cells.Item(97, 8).Value = 4276932  # H97: 4447981.5 -> 4276932
$ws.Cells.Item(97, 10).Value = 15877400  # J97: 18523518 -> 15877400
$ws.Cells.Item(97, 12).Value = 15877400  # L97: 18523518 -> 15877400
$ws.Cells.Item(97, 14).Value = -15878392  # N97: -18524510 -> -15878392
$ws.Cells.Item(122, 8).Value = 5269.067  # H122: 5257.9556 -> 5269.067
$ws.Cells.Item(122, 9).Value = 5142.5264  # I122: 5168.8423 -> 5142.5264
$ws.Cells.Item(122, 10).Value = 5361.5386  # J122: 5323.077 -> 5361.5386
$ws.Cells.Item(122, 11).Value = 15427.5792  # K122: 15506.5269 -> 15427.5792
$ws.Cells.Item(122, 12).Value = 16084.6158  # L122: 15969.231 -> 16084.6158
$ws.Cells.Item(122, 13).Value = -12977.5792  # M122: -13056.5269 -> -12977.5792
$ws.Cells.Item(122, 14).Value = -20984.6158  # N122: -20869.231 -> -20984.6158
$ws.Cells.Item(125, 8).Value = 29833.334  # H125: 29769.23 -> 29833.334
$ws.Cells.Item(125, 10).Value = 29833.334  # J125: 29769.23 -> 29833.334
$ws.Cells.Item(125, 12).Value = 29833.334  # L125: 29769.23 -> 29833.334
$ws.Cells.Item(125, 14).Value = -39673.334  # N125: -39609.23 -> -39673.334
$ws.Cells.Item(132, 8).Value = 2553  # H132: 2250 -> 2553
$ws.Cells.Item(132, 9).Value = 2553  # I132: 2250 -> 2553
$ws.Cells.Item(132, 11).Value = 7659  # K132: 6750 -> 7659
$ws.Cells.Item(132, 13).Value = -5129  # M132: -4220 -> -5129
$ws.Cells.Item(136, 8).Value = 4571.433  # H136: 4010.3845 -> 4571.433
$ws.Cells.Item(136, 9).Value = 1922.25  # I136: 1843.3334 -> 1922.25
$ws.Cells.Item(136, 10).Value = 5534.773  # J136: 8886.25 -> 5534.773
$ws.Cells.Item(136, 11).Value = 5766.75  # K136: 5530.0002 -> 5766.75
$ws.Cells.Item(136, 12).Value = 16604.319  # L136: 26658.75 -> 16604.319
$ws.Cells.Item(136, 13).Value = -3216.75  # M136: -2980.0002 -> -3216.75
$ws.Cells.Item(136, 14).Value = -21704.319  # N136: -31758.75 -> -21704.319

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(62, 8).Value = 50000  # H62: 0 -> 50000
$ws.Cells.Item(62, 10).Value = 50000  # J62: 0 -> 50000
$ws.Cells.Item(62, 12).Value = 50000  # L62: 0 -> 50000
$ws.Cells.Item(62, 14).Value = -51372  # N62: None -> -51372
$ws.Cells.Item(65, 8).Value = 50000  # H65: 0 -> 50000
$ws.Cells.Item(65, 10).Value = 50000  # J65: 0 -> 50000
$ws.Cells.Item(65, 12).Value = 150000  # L65: 0 -> 150000
$ws.Cells.Item(65, 14).Value = -156864  # N65: None -> -156864
$ws.Cells.Item(94, 8).Value = 2648331.2  # H94: 2586875.8 -> 2648331.2
$ws.Cells.Item(94, 9).Value = 2255.862  # I94: 2343.8276 -> 2255.862
$ws.Cells.Item(94, 10).Value = 8551115  # J94: 7940549 -> 8551115
$ws.Cells.Item(94, 11).Value = 2255.862  # K94: 2343.8276 -> 2255.862
$ws.Cells.Item(94, 12).Value = 8551115  # L94: 7940549 -> 8551115
$ws.Cells.Item(94, 13).Value = -1804.862  # M94: -1892.8276 -> -1804.862
$ws.Cells.Item(94, 14).Value = -8552017  # N94: -7941451 -> -8552017
$ws.Cells.Item(107, 8).Value = 13928  # H107: 10020.333 -> 13928
$ws.Cells.Item(107, 9).Value = 33755.5  # I107: 13787.143 -> 33755.5
$ws.Cells.Item(107, 10).Value = 7318.8335  # J107: 6724.375 -> 7318.8335
$ws.Cells.Item(107, 11).Value = 33755.5  # K107: 13787.143 -> 33755.5
$ws.Cells.Item(107, 12).Value = 7318.8335  # L107: 6724.375 -> 7318.8335
$ws.Cells.Item(107, 13).Value = -31835.5  # M107: -11867.143 -> -31835.5
$ws.Cells.Item(107, 14).Value = -11158.8335  # N107: -10564.375 -> -11158.8335
$ws.Cells.Item(115, 8).Value = 20000  # H115: 0 -> 20000
$ws.Cells.Item(115, 10).Value = 20000  # J115: 0 -> 20000
$ws.Cells.Item(115, 12).Value = 20000  # L115: 0 -> 20000
$ws.Cells.Item(115, 14).Value = -23134  # N115: None -> -23134
$ws.Cells.Item(134, 8).Value = 4197.625  # H134: 5304 -> 4197.625
$ws.Cells.Item(134, 9).Value = 4197.625  # I134: 5304 -> 4197.625
$ws.Cells.Item(134, 11).Value = 12592.875  # K134: 15912 -> 12592.875
$ws.Cells.Item(134, 13).Value = -10057.875  # M134: -13377 -> -10057.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1346.4  # H94: 1385.6428 -> 1346.4
$ws.Cells.Item(94, 9).Value = 1017.5  # I94: 1061.6 -> 1017.5
$ws.Cells.Item(94, 11).Value = 1017.5  # K94: 1061.6 -> 1017.5
$ws.Cells.Item(94, 13).Value = -566.5  # M94: -610.5999999999999 -> -566.5
$ws.Cells.Item(99, 8).Value = 2749.5  # H99: 0 -> 2749.5
$ws.Cells.Item(99, 9).Value = 2749.5  # I99: 0 -> 2749.5
$ws.Cells.Item(99, 11).Value = 2749.5  # K99: 0 -> 2749.5
$ws.Cells.Item(99, 13).Value = -1251.5  # M99: None -> -1251.5
$ws.Cells.Item(122, 8).Value = 4302.579  # H122: 4608.7646 -> 4302.579
$ws.Cells.Item(122, 9).Value = 3841.7856  # I122: 4198.75 -> 3841.7856
$ws.Cells.Item(122, 11).Value = 11525.3568  # K122: 12596.25 -> 11525.3568
$ws.Cells.Item(122, 13).Value = -9075.356800000001  # M122: -10146.25 -> -9075.356800000001
$ws.Cells.Item(126, 8).Value = 2749.5  # H126: 0 -> 2749.5
$ws.Cells.Item(126, 9).Value = 2749.5  # I126: 0 -> 2749.5
$ws.Cells.Item(126, 11).Value = 8248.5  # K126: 0 -> 8248.5
$ws.Cells.Item(126, 13).Value = -5778.5  # M126: None -> -5778.5
$ws.Cells.Item(132, 8).Value = 6000  # H132: 3190.6667 -> 6000
$ws.Cells.Item(132, 9).Value = 6000  # I132: 3208 -> 6000
$ws.Cells.Item(132, 10).Value = 0  # J132: 3000 -> 0
$ws.Cells.Item(132, 11).Value = 18000  # K132: 9624 -> 18000
$ws.Cells.Item(132, 12).Value = 0  # L132: 9000 -> 0
$ws.Cells.Item(132, 13).Value = -15470  # M132: -7094 -> -15470
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -14060

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1069.9166  # H34: 949.0909 -> 1069.9166
$ws.Cells.Item(34, 10).Value = 1999.8334  # J34: 1920 -> 1999.8334
$ws.Cells.Item(34, 12).Value = 5999.5002  # L34: 5760 -> 5999.5002
$ws.Cells.Item(34, 14).Value = -6167.5002  # N34: -5928 -> -6167.5002
$ws.Cells.Item(37, 8).Value = 1099724.5  # H37: 793148.7 -> 1099724.5
$ws.Cells.Item(37, 10).Value = 1099724.5  # J37: 793148.7 -> 1099724.5
$ws.Cells.Item(37, 12).Value = 3299173.5  # L37: 2379446.1 -> 3299173.5
$ws.Cells.Item(37, 14).Value = -3299397.5  # N37: -2379670.1 -> -3299397.5
$ws.Cells.Item(45, 8).Value = 1727  # H45: 1849.2222 -> 1727
$ws.Cells.Item(45, 9).Value = 1804.25  # I45: 1961.75 -> 1804.25
$ws.Cells.Item(45, 10).Value = 1649.75  # J45: 1759.2 -> 1649.75
$ws.Cells.Item(45, 11).Value = 5412.75  # K45: 5885.25 -> 5412.75
$ws.Cells.Item(45, 12).Value = 4949.25  # L45: 5277.6 -> 4949.25
$ws.Cells.Item(45, 13).Value = -4880.75  # M45: -5353.25 -> -4880.75
$ws.Cells.Item(45, 14).Value = -6013.25  # N45: -6341.6 -> -6013.25
$ws.Cells.Item(46, 8).Value = 200530  # H46: 250412.5 -> 200530
$ws.Cells.Item(46, 10).Value = 1000  # J46: 0 -> 1000
$ws.Cells.Item(46, 12).Value = 3000  # L46: 0 -> 3000
$ws.Cells.Item(46, 14).Value = -3182  # N46: None -> -3182

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 111484820  # H80: 83616110 -> 111484820
$ws.Cells.Item(80, 9).Value = 1111111  # I80: 560555 -> 1111111
$ws.Cells.Item(80, 11).Value = 1111111  # K80: 560555 -> 1111111
$ws.Cells.Item(80, 13).Value = -1110113  # M80: -559557 -> -1110113
$ws.Cells.Item(83, 8).Value = 111484820  # H83: 83616110 -> 111484820
$ws.Cells.Item(83, 9).Value = 1111111  # I83: 560555 -> 1111111
$ws.Cells.Item(83, 11).Value = 5555555  # K83: 2802775 -> 5555555
$ws.Cells.Item(83, 13).Value = -5550563  # M83: -2797783 -> -5550563
$ws.Cells.Item(122, 8).Value = 4239.7085  # H122: 4435.857 -> 4239.7085
$ws.Cells.Item(122, 9).Value = 3645.6924  # I122: 3879.4 -> 3645.6924
$ws.Cells.Item(122, 11).Value = 10937.0772  # K122: 11638.2 -> 10937.0772
$ws.Cells.Item(122, 13).Value = -8487.0772  # M122: -9188.200000000001 -> -8487.0772
$ws.Cells.Item(126, 8).Value = 7562.55  # H126: 7786.9473 -> 7562.55
$ws.Cells.Item(126, 9).Value = 5764.7144  # I126: 6175.6665 -> 5764.7144
$ws.Cells.Item(126, 11).Value = 17294.1432  # K126: 18526.9995 -> 17294.1432
$ws.Cells.Item(126, 13).Value = -14824.1432  # M126: -16056.9995 -> -14824.1432
$ws.Cells.Item(132, 8).Value = 6548.205  # H132: 6913.222 -> 6548.205
$ws.Cells.Item(132, 9).Value = 6364  # I132: 6630.8125 -> 6364
$ws.Cells.Item(132, 10).Value = 7800.8  # J132: 9172.5 -> 7800.8
$ws.Cells.Item(132, 11).Value = 19092  # K132: 19892.4375 -> 19092
$ws.Cells.Item(132, 12).Value = 23402.4  # L132: 27517.5 -> 23402.4
$ws.Cells.Item(132, 13).Value = -16562  # M132: -17362.4375 -> -16562
$ws.Cells.Item(132, 14).Value = -28462.4  # N132: -32577.5 -> -28462.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2923.5454  # H7: 3125.8438 -> 2923.5454
$ws.Cells.Item(7, 9).Value = 2970.5  # I7: 3139.7307 -> 2970.5
$ws.Cells.Item(7, 10).Value = 2749.1428  # J7: 3065.6667 -> 2749.1428
$ws.Cells.Item(7, 11).Value = 2970.5  # K7: 3139.7307 -> 2970.5
$ws.Cells.Item(7, 12).Value = 2749.1428  # L7: 3065.6667 -> 2749.1428
$ws.Cells.Item(7, 13).Value = -2858.5  # M7: -3027.7307 -> -2858.5
$ws.Cells.Item(7, 14).Value = -2973.1428  # N7: -3289.6667 -> -2973.1428
$ws.Cells.Item(40, 8).Value = 4258.8945  # H40: 4385 -> 4258.8945
$ws.Cells.Item(40, 9).Value = 4032.7856  # I40: 4190 -> 4032.7856
$ws.Cells.Item(40, 11).Value = 4032.7856  # K40: 4190 -> 4032.7856
$ws.Cells.Item(40, 13).Value = -3896.7856  # M40: -4054 -> -3896.7856
$ws.Cells.Item(93, 8).Value = 5052474  # H93: 5557686 -> 5052474
$ws.Cells.Item(93, 9).Value = 1723.6897  # I93: 1881.5769 -> 1723.6897
$ws.Cells.Item(93, 11).Value = 1723.6897  # K93: 1881.5769 -> 1723.6897
$ws.Cells.Item(93, 13).Value = -475.6896999999999  # M93: -633.5769 -> -475.6896999999999
$ws.Cells.Item(100, 8).Value = 213156.5  # H100: 194051.19 -> 213156.5
$ws.Cells.Item(100, 9).Value = 213156.5  # I100: 194051.19 -> 213156.5
$ws.Cells.Item(100, 11).Value = 213156.5  # K100: 194051.19 -> 213156.5
$ws.Cells.Item(100, 13).Value = -212615.5  # M100: -193510.19 -> -212615.5
$ws.Cells.Item(126, 8).Value = 2923.5454  # H126: 3125.8438 -> 2923.5454
$ws.Cells.Item(126, 9).Value = 2970.5  # I126: 3139.7307 -> 2970.5
$ws.Cells.Item(126, 10).Value = 2749.1428  # J126: 3065.6667 -> 2749.1428
$ws.Cells.Item(126, 11).Value = 8911.5  # K126: 9419.1921 -> 8911.5
$ws.Cells.Item(126, 12).Value = 8247.428400000001  # L126: 9197.000100000001 -> 8247.428400000001
$ws.Cells.Item(126, 13).Value = -6441.5  # M126: -6949.1921 -> -6441.5
$ws.Cells.Item(126, 14).Value = -13187.4284  # N126: -14137.0001 -> -13187.4284
$ws.Cells.Item(132, 8).Value = 6791.0386  # H132: 7541.0967 -> 6791.0386
$ws.Cells.Item(132, 9).Value = 6398.5  # I132: 7256.8423 -> 6398.5
$ws.Cells.Item(132, 10).Value = 7674.25  # J132: 7991.1665 -> 7674.25
$ws.Cells.Item(132, 11).Value = 19195.5  # K132: 21770.5269 -> 19195.5
$ws.Cells.Item(132, 12).Value = 23022.75  # L132: 23973.4995 -> 23022.75
$ws.Cells.Item(132, 13).Value = -16665.5  # M132: -19240.5269 -> -16665.5
$ws.Cells.Item(132, 14).Value = -28082.75  # N132: -29033.4995 -> -28082.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(86, 8).Value = 80108.164  # H86: 86912.5 -> 80108.164
$ws.Cells.Item(86, 10).Value = 80108.164  # J86: 86912.5 -> 80108.164
$ws.Cells.Item(86, 12).Value = 80108.164  # L86: 86912.5 -> 80108.164
$ws.Cells.Item(86, 14).Value = -82354.164  # N86: -89158.5 -> -82354.164
$ws.Cells.Item(87, 8).Value = 20000  # H87: 40000 -> 20000
$ws.Cells.Item(87, 10).Value = 20000  # J87: 40000 -> 20000
$ws.Cells.Item(87, 12).Value = 20000  # L87: 40000 -> 20000
$ws.Cells.Item(87, 14).Value = -22496  # N87: -42496 -> -22496
$ws.Cells.Item(89, 8).Value = 80108.164  # H89: 86912.5 -> 80108.164
$ws.Cells.Item(89, 10).Value = 80108.164  # J89: 86912.5 -> 80108.164
$ws.Cells.Item(89, 12).Value = 400540.82  # L89: 434562.5 -> 400540.82
$ws.Cells.Item(89, 14).Value = -411772.82  # N89: -445794.5 -> -411772.82
$ws.Cells.Item(90, 8).Value = 20000  # H90: 40000 -> 20000
$ws.Cells.Item(90, 10).Value = 20000  # J90: 40000 -> 20000
$ws.Cells.Item(90, 12).Value = 60000  # L90: 120000 -> 60000
$ws.Cells.Item(90, 14).Value = -72480  # N90: -132480 -> -72480
$ws.Cells.Item(93, 8).Value = 30000  # H93: 29333.334 -> 30000
$ws.Cells.Item(93, 10).Value = 30000  # J93: 29333.334 -> 30000
$ws.Cells.Item(93, 12).Value = 30000  # L93: 29333.334 -> 30000
$ws.Cells.Item(93, 14).Value = -34992  # N93: -34325.334 -> -34992
$ws.Cells.Item(122, 8).Value = 1480.8182  # H122: 1509.5 -> 1480.8182
$ws.Cells.Item(122, 9).Value = 1398.4286  # I122: 1432.5 -> 1398.4286
$ws.Cells.Item(122, 11).Value = 4195.2858  # K122: 4297.5 -> 4195.2858
$ws.Cells.Item(122, 13).Value = -1745.2858  # M122: -1847.5 -> -1745.2858
